$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "VENTAS POR GRUPO" - update new sale amounts for row 19
#    (MANCHENO PINO HERVIN SANTIAGO) and "de 31" counters on row 33
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H19").Value = 71.09999999999999
$wsGrupo.Range("I19").Value = 70.2
$wsGrupo.Range("L19").Value = 221.72
$wsGrupo.Range("M19").Value = 42.04
$wsGrupo.Range("H33").Value = "1 de 31"
$wsGrupo.Range("I33").Value = "1 de 31"
$wsGrupo.Range("L33").Value = "1 de 31"
$wsGrupo.Range("M33").Value = "1 de 31"

# ---------------------------------------------------------------
# 2) "VENTA MENSUAL" - record the new "agosto" sale for the same
#    client and refresh the column total
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F19").Value = 405.06
$wsMensual.Range("F33").Value = 862.98

# ---------------------------------------------------------------
# 3) "CUMPLIMIENTO MENSUAL" - a new GRUPO ("240X120 PORCELANATO")
#    is now reported; insert its row (pushing the rest down by one)
#    and refresh every PRESUPUESTO/VENTA/POR CUMPLIR/CUMPLIMIENTO
#    figure for the advisor to match the latest recalculated report.
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Insert the new row 2 and clone the number formatting of the row
# that lands right below it so currency/percentage formats match.
$wsCumpl.Rows.Item(2).Insert()
$wsCumpl.Range("A3:F3").Copy()
$wsCumpl.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: 240X120 PORCELANATO
$wsCumpl.Range("A2").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B2").Value = "240X120 PORCELANATO"
$wsCumpl.Range("C2").Value = 260.285000070615
$wsCumpl.Range("D2").Value = 0
$wsCumpl.Range("E2").Value = 260.285000070615
$wsCumpl.Range("F2").Value = 0

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("A3").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B3").Value = "240X80 PORCELANATO"
$wsCumpl.Range("C3").Value = 3120.1145
$wsCumpl.Range("D3").Value = 457.92
$wsCumpl.Range("E3").Value = 2662.1945
$wsCumpl.Range("F3").Value = 0.1467638447242882

# Row 4: FREGADEROS DE COCINA
$wsCumpl.Range("A4").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B4").Value = "FREGADEROS DE COCINA"
$wsCumpl.Range("C4").Value = 646.361575487259
$wsCumpl.Range("D4").Value = 0
$wsCumpl.Range("E4").Value = 646.361575487259
$wsCumpl.Range("F4").Value = 0

# Row 5: GRANITO
$wsCumpl.Range("A5").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B5").Value = "GRANITO"
$wsCumpl.Range("C5").Value = 238.32
$wsCumpl.Range("D5").Value = 0
$wsCumpl.Range("E5").Value = 238.32
$wsCumpl.Range("F5").Value = 0

# Row 6: GRIFERIAS
$wsCumpl.Range("A6").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B6").Value = "GRIFERIAS"
$wsCumpl.Range("C6").Value = 106.82
$wsCumpl.Range("D6").Value = 0
$wsCumpl.Range("E6").Value = 106.82
$wsCumpl.Range("F6").Value = 0

# Row 7: INODOROS
$wsCumpl.Range("A7").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B7").Value = "INODOROS"
$wsCumpl.Range("C7").Value = 1600
$wsCumpl.Range("D7").Value = 71.09999999999999
$wsCumpl.Range("E7").Value = 1528.9
$wsCumpl.Range("F7").Value = 0.0444375

# Row 8: LAVABOS
$wsCumpl.Range("A8").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B8").Value = "LAVABOS"
$wsCumpl.Range("C8").Value = 625
$wsCumpl.Range("D8").Value = 70.2
$wsCumpl.Range("E8").Value = 554.8
$wsCumpl.Range("F8").Value = 0.11232

# Row 9: LED
$wsCumpl.Range("A9").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B9").Value = "LED"
$wsCumpl.Range("C9").Value = 300
$wsCumpl.Range("D9").Value = 0
$wsCumpl.Range("E9").Value = 300
$wsCumpl.Range("F9").Value = 0

# Row 10: NO RESURTIBLES
$wsCumpl.Range("A10").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B10").Value = "NO RESURTIBLES"
$wsCumpl.Range("C10").Value = 650.25
$wsCumpl.Range("D10").Value = 0
$wsCumpl.Range("E10").Value = 650.25
$wsCumpl.Range("F10").Value = 0

# Row 11: OTROS
$wsCumpl.Range("A11").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B11").Value = "OTROS"
$wsCumpl.Range("C11").Value = 0
$wsCumpl.Range("D11").Value = 0
$wsCumpl.Range("E11").Value = 0
$wsCumpl.Range("F11").Value = 0

# Row 12: PANELES DECORATIVOS
$wsCumpl.Range("A12").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B12").Value = "PANELES DECORATIVOS"
$wsCumpl.Range("C12").Value = 100
$wsCumpl.Range("D12").Value = 0
$wsCumpl.Range("E12").Value = 100
$wsCumpl.Range("F12").Value = 0

# Row 13: PANELES PU
$wsCumpl.Range("A13").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B13").Value = "PANELES PU"
$wsCumpl.Range("C13").Value = 20
$wsCumpl.Range("D13").Value = 0
$wsCumpl.Range("E13").Value = 20
$wsCumpl.Range("F13").Value = 0

# Row 14: PANELES PVC
$wsCumpl.Range("A14").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B14").Value = "PANELES PVC"
$wsCumpl.Range("C14").Value = 100
$wsCumpl.Range("D14").Value = 0
$wsCumpl.Range("E14").Value = 100
$wsCumpl.Range("F14").Value = 0

# Row 15: PIEDRA SINTERIZADA
$wsCumpl.Range("A15").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B15").Value = "PIEDRA SINTERIZADA"
$wsCumpl.Range("C15").Value = 527.03
$wsCumpl.Range("D15").Value = 221.72
$wsCumpl.Range("E15").Value = 305.3099999999999
$wsCumpl.Range("F15").Value = 0.4206971140162799

# Row 16: PORCELANATO
$wsCumpl.Range("A16").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B16").Value = "PORCELANATO"
$wsCumpl.Range("C16").Value = 21873.1
$wsCumpl.Range("D16").Value = 42.04
$wsCumpl.Range("E16").Value = 21831.06
$wsCumpl.Range("F16").Value = 0.001921995510467195

# Row 17: PUERTAS DE SEGURIDAD
$wsCumpl.Range("A17").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B17").Value = "PUERTAS DE SEGURIDAD"
$wsCumpl.Range("C17").Value = 342
$wsCumpl.Range("D17").Value = 0
$wsCumpl.Range("E17").Value = 342
$wsCumpl.Range("F17").Value = 0

# Row 18: SAL SOLUBLE
$wsCumpl.Range("A18").Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$wsCumpl.Range("B18").Value = "SAL SOLUBLE"
$wsCumpl.Range("C18").Value = 1600
$wsCumpl.Range("D18").Value = 0
$wsCumpl.Range("E18").Value = 1600
$wsCumpl.Range("F18").Value = 0

# Row 19: TOTAL
$wsCumpl.Range("B19").Value = "TOTAL"
$wsCumpl.Range("C19").Value = 32109.28107555787
$wsCumpl.Range("D19").Value = 862.98
$wsCumpl.Range("E19").Value = 31246.30107555788
$wsCumpl.Range("F19").Value = 0.02687634139080476
